$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.99999999030758002
$ws.Range("A2").Value = 0.99622868625103644
$ws.Range("A3").Value = 0.97968973441769225
$ws.Range("A4").Value = 0.97333673064905679
$ws.Range("A5").Value = 0.96400335616536159
$ws.Range("A6").Value = 0.94129960412871738
$ws.Range("A7").Value = 0.93764679048191768
$ws.Range("A8").Value = 0.93329213564236979
$ws.Range("A9").Value = 0.92959149997415136
$ws.Range("A10").Value = 0.92229772598362281
$ws.Range("A11").Value = 0.92279485442630138
$ws.Range("A12").Value = 0.92399864508111818
$ws.Range("A13").Value = 0.9127117542771177
$ws.Range("A14").Value = 0.9085445003015975
$ws.Range("A15").Value = 0.90595307484802112
$ws.Range("A16").Value = 0.90344661141747173
$ws.Range("A17").Value = 0.89973873237561497
$ws.Range("A18").Value = 0.89862983577890054
$ws.Range("A19").Value = 0.9951755107544159
$ws.Range("A20").Value = 0.98805862089253083
$ws.Range("A21").Value = 0.98666015717935007
$ws.Range("A22").Value = 0.98539565542986507
$ws.Range("A23").Value = 0.9726546734870225
$ws.Range("A24").Value = 0.95963351728773061
$ws.Range("A25").Value = 0.95317648725796433
$ws.Range("A26").Value = 0.93873188398700913
$ws.Range("A27").Value = 0.93540133904819922
$ws.Range("A28").Value = 0.92093006086901208
$ws.Range("A29").Value = 0.91085549700147539
$ws.Range("A30").Value = 0.90688389986322426
$ws.Range("A31").Value = 0.90582623037124188
$ws.Range("A32").Value = 0.90116138300419946
$ws.Range("A33").Value = 0.90064139148267519
